# Update the Price (column D) and Volume(1h) (column E) values for the
# cryptos sheet, per the commit "Updated symbol list on Fri Jan 20 23:35:57
# UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = "..."; E = "..." }  (only keys present are changed)
$changes = @{
    2  = @{ D = "302.96";         E = "2.83%" }
    3  = @{ D = "33.62";          E = "8.33%" }
    4  = @{ D = "5.099";          E = "3.41%" }
    5  = @{ D = "0.07802";        E = "5.92%" }
    6  = @{ D = "2.408";          E = "4.36%" }
    7  = @{ D = "8.022";          E = "3.98%" }
    8  = @{ D = "3.908";          E = "3.95%" }
    9  = @{ D = "0.9361";         E = "2.54%" }
    10 = @{ D = "0.09820";        E = "18.53%" }
    11 = @{ D = "0.1781";         E = "5.64%" }
    12 = @{ D = "0.08550";        E = "3.58%" }
    13 = @{ D = "0.03340";        E = "7.07%" }
    14 = @{ D = "0.09901";        E = "-1.69%" }
    15 = @{ D = "0.001475";       E = "-2.16%" }
    16 = @{ D = "0.005805";       E = "1.11%" }
    17 = @{ D = "3.464";          E = "-0.53%" }
    18 = @{ D = "2.187";          E = "5.17%" }
    19 = @{ D = "0.3368";         E = "1.16%" }
    20 = @{ D = "0.1347";         E = "3.35%" }
    21 = @{ D = "4.295";          E = "7.99%" }
    22 = @{ D = "0.2286";         E = "8.89%" }
    23 = @{ D = "0.04644";        E = "1.99%" }
    24 = @{ D = "0.001217";       E = "0.59%" }
    25 = @{ D = "0.004414";       E = "1.68%" }
    26 = @{ D = "0.0001293";      E = "-0.48%" }
    27 = @{ E = "-0.46%" }
    39 = @{ D = "0.01741";        E = "8.26%" }
    40 = @{ D = "0.04814";        E = "8.31%" }
    41 = @{ D = "0.007751";       E = "5.67%" }
    42 = @{ D = "0.009720";       E = "10.01%" }
    43 = @{ D = "0.1409";         E = "6.31%" }
    44 = @{ D = "0.002067";       E = "0.40%" }
    45 = @{ D = "0.009092";       E = "-0.11%" }
    46 = @{ E = "0.63%" }
    47 = @{ D = "0.00000000746"; E = "-0.45%" }
    48 = @{ D = "2.551";          E = "13.84%" }
    49 = @{ E = "-31.35%" }
    50 = @{ D = "0.00002090";     E = "-0.45%" }
    51 = @{ D = "0.0001990";      E = "-0.45%" }
}

# Values such as "302.96" or "2.83%" look numeric/percentage to Excel, so a
# plain .Value assignment would silently convert them to a number (losing the
# exact textual representation, e.g. trailing zeros, "%" sign). Forcing the
# cell to Text format before the assignment keeps the literal string, and
# resetting the style back to Normal afterwards avoids leaving a stray
# number-format style applied to the cell (matching the original formatting).
foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
        $cell.Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
        $cell.Style = "Normal"
    }
}
